# Connected Office Test Data - mark Create/Update (and for device rows,
# Delete) tests as Passed ("TRUE"), and append a row for a newly-run
# Update Test / Create Test.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Results")

$xlPasteValues = -4163

# Rows 2-24 hold the Zone (2-10), Category (11-15) and Device (16-24)
# test rows. Column C ("Read Test Passed") already reads "TRUE" as text
# for every row, so re-use it as the source of a literal "TRUE" string
# (assigning the word "TRUE" straight to .Value would be auto-typed as a
# Boolean instead of text).
$sourceCell = $ws.Range("C2")
$sourceCell.Copy()

for ($row = 2; $row -le 24; $row++) {
    $ws.Range("B" + $row).PasteSpecial($xlPasteValues)
    $ws.Range("D" + $row).PasteSpecial($xlPasteValues)

    # The Device rows (16-24) also had their Delete Test still failing;
    # mark those passed too.
    if ($row -ge 16) {
        $ws.Range("E" + $row).PasteSpecial($xlPasteValues)
    }
}

$excel.CutCopyMode = 0

# Record the results of the newly added "Update Test" / "Create Test"
# run as a fresh row underneath the existing data.
$newRow = $ws.Range("B25")
$newRow.Style = "Normal"
$sourceCell.Copy()
$newRow.PasteSpecial($xlPasteValues)

$excel.CutCopyMode = 0
